$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New skill: "Merchants Coin"
$skillName = "Merchants Coin"
$skillDesc = "The one good thing about being a merchant is defending your coin no matter what it takes. The other good thing, is that people you need to fight beside you will help you. Move quicker and fight faster!"

# Row 27
$ws.Range("A27").Value = 27
$ws.Range("B27").Value = $skillName
$ws.Range("D27").Value = $skillDesc
$ws.Range("E27").Value = 999
$ws.Range("F27").Value = 0.015
$ws.Range("I27").Value = 0.03
$ws.Range("J27").Value = 0.05
$ws.Range("K27").Value = 0.0001
$ws.Range("O27").Value = 1
$ws.Range("Q27").Value = 11
$ws.Range("R27").Value = 0

# Row 28
$ws.Range("A28").Value = 28
$ws.Range("B28").Value = $skillName
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = $skillDesc
$ws.Range("E28").Value = 999
$ws.Range("F28").Value = 0.00095
$ws.Range("I28").Value = 0.0003
$ws.Range("J28").Value = 0.0005
$ws.Range("K28").Value = 0.0001
$ws.Range("O28").Value = 1
$ws.Range("Q28").Value = 11
$ws.Range("R28").Value = 0
